$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing rows 145+ down by one.
$ws.Rows.Item(145).Insert()

# Copy the formatting (white "Proposed" fill) from an existing "Proposed"-styled
# row (row 73) onto the newly inserted row 145, so it reuses the same style
# instead of creating a brand new one.
$ws.Range("A73:V73").Copy()
$ws.Range("A145:V145").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row with the "post-traumatic stress symptom severity scale"
# entry.
$ws.Range("A145").Value = "GMHO:0000263"
$ws.Range("B145").Value = "post-traumatic stress symptom severity scale"
$ws.Range("C145").Value = "A measurement scale that is used to measure post-traumatic stress symptom severity."
$ws.Range("D145").Value = "measurement scale"
$ws.Range("J145").Value = "Intervention outcomes and spillover effects"
$ws.Range("P145").Value = "LSR 2"
$ws.Range("Q145").Value = "Intervention outcomes and spillover effects"
$ws.Range("S145").Value = "Proposed"
